$d = $word.ActiveDocument

# --- 1. Remove the "Meta description" paragraph that currently sits right
#        after the H1 title at the top of the document. ---
$metaPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Meta description:*") {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# --- 2. Insert a new bold "Play Free Auspicious Fortune God Slot" paragraph
#        right before the final paragraph, and turn the final paragraph's
#        text into the meta-description text (keeping its italic run). ---
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)

$lastPara.Range.InsertParagraphBefore()

$n2 = $d.Paragraphs.Count
$titlePara = $d.Paragraphs.Item($n2 - 1)
$titlePara.Range.Text = "Play Free Auspicious Fortune God Slot"
$titlePara.Range.Font.Bold = 1

$finalPara = $d.Paragraphs.Item($n2)
$finalRange = $finalPara.Range
$find = $finalRange.Find
$find.Execute("Create an eye-catching feature image for " + [char]34 + "Auspicious Fortune God" + [char]34 + " that showcases a happy Maya warrior with glasses in cartoon style. Use bright and bold colors that are reminiscent of the game's Chinese-theme, such as red and gold, to catch the viewer's attention. Position the warrior in a confident and charismatic pose, with his hands on his hips and a big smile on his face. Make sure to incorporate the game's logo into the image and any other relevant symbols such as fortune deities and money trees. The overall image should convey a sense of excitement and fun while highlighting the game's unique features.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Read our review of Auspicious Fortune God, a Chinese-themed slot game with Expanding Wilds, Respins, and Customization Options. Play for free.", `
    2)
